$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in F1:H1, matching style of existing headers (A1:E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Boolean values per row for columns F (KNN), G (SVM), H (RF)
$values = @(
    @($false, $false, $false),  # row 2
    @($true,  $false, $false),  # row 3
    @($true,  $false, $false),  # row 4
    @($false, $false, $false),  # row 5
    @($false, $false, $false),  # row 6
    @($false, $false, $false),  # row 7
    @($false, $false, $false),  # row 8
    @($false, $false, $false),  # row 9
    @($false, $false, $false),  # row 10
    @($false, $false, $false),  # row 11
    @($false, $false, $false)   # row 12
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i][0]
    $ws.Cells.Item($row, 7).Value = $values[$i][1]
    $ws.Cells.Item($row, 8).Value = $values[$i][2]
}
